$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 797
$ws.Range("I2").Value = 2216
$ws.Range("J2").Value = 9243
$ws.Range("K2").Value = 44
$ws.Range("L2").Value = 2551
$ws.Range("M2").Value = 152
$ws.Range("N2").Value = 1704
$ws.Range("O2").Value = 5
$ws.Range("P2").Value = 41
$ws.Range("Q2").Value = 17
$ws.Range("R2").Value = 144
$ws.Range("S2").Value = 995
$ws.Range("T2").Value = 1589
$ws.Range("U2").Value = 135
$ws.Range("V2").Value = 14435
$ws.Range("W2").Value = 8
$ws.Range("X2").Value = 14546
$ws.Range("Y2").Value = 15
$ws.Range("Z2").Value = 228
$ws.Range("AA2").Value = 105
